# Updated cryptos list on Sat Mar  2 19:50:00 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns of the cryptos table with
# new quotes, and fixes the ARBITRUM / WEMIXToken rows (45-46), which were
# swapped, along with their refreshed price/volume figures.
#
# Price-looking values that would otherwise be auto-detected as numbers by
# Excel (e.g. "1.00", "130.09") are forced back to Text via NumberFormat
# "@" so the stored cell content matches the source data exactly (values
# like "62.068.77" that contain more than one '.' are already left alone
# as text by Excel's auto-detection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.068.77'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '3.428.36'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.09'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  +6.46%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.739'
$ws.Range("E9").Value = '  +6.98%  '
$ws.Range("E10").Value = '  +4.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.92'
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000225'
$ws.Range("E12").Value = '  +50.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.22'
$ws.Range("E13").Value = '  +10.06%  '
$ws.Range("D15").Value = '3.976.75'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("E16").Value = '  +7.69%  '
$ws.Range("D17").Value = '3.423.93'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.54'
$ws.Range("E18").Value = '  +8.07%  '
$ws.Range("E19").Value = '  +7.65%  '
$ws.Range("D20").Value = '62.034.33'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '456.08'
$ws.Range("E21").Value = '  +46.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.59'
$ws.Range("E22").Value = '  +8.97%  '
$ws.Range("E23").Value = '  +1.41%  '
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.29'
$ws.Range("E25").Value = '  +3.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.17'
$ws.Range("E26").Value = '  +11.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.13'
$ws.Range("E27").Value = '  +11.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.78'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.65'
$ws.Range("E29").Value = '  -2.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.13'
$ws.Range("E30").Value = '  +6.72%  '
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.05'
$ws.Range("E33").Value = '  -3.44%  '
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  +3.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.30'
$ws.Range("E37").Value = '  +4.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("E39").Value = '  +1.89%  '
$ws.Range("E40").Value = '  +7.86%  '
$ws.Range("E41").Value = '  -1.24%  '
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.44'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  +9.60%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.55'
$ws.Range("E45").Value = '  +15.22%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.00'
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.65'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.48'
$ws.Range("E48").Value = '  +6.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.142'
$ws.Range("E49").Value = '  +19.30%  '
$ws.Range("D50").Value = '3.777.36'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.11'
$ws.Range("E51").Value = '  +7.51%  '
